$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on columns that must stay text-typed even though
# their content looks numeric / date-like, matching the inlineStr cells
# added in the diff (row 22). Revert the cell style afterwards so the
# only lasting effect is the text storage type, not a visual format change.
foreach ($addr in @("I22", "Y22", "Z22", "AA22", "AB22")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A22").Value = 112183052
$ws.Range("B22").Value = 77597
$ws.Range("C22").Value = "Ovaliderad"
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 864
$ws.Range("F22").Value = "Knottrig blåslav"
$ws.Range("G22").Value = "Hypogymnia bitteri"
$ws.Range("H22").Value = "(Lynge) Ahti"
$ws.Range("I22").Value = "1"
$ws.Range("P22").Value = "Tjappsåive, Pi lm"
$ws.Range("Q22").Value = 699856.6862899091
$ws.Range("R22").Value = 7309603.113012934
$ws.Range("S22").Value = 5
$ws.Range("T22").Value = "Norrbotten"
$ws.Range("U22").Value = "Arvidsjaur"
$ws.Range("V22").Value = "Pite lappmark"
$ws.Range("W22").Value = "Arvidsjaur"
$ws.Range("Y22").Value = "2023-08-17"
$ws.Range("Z22").Value = "00:00"
$ws.Range("AA22").Value = "2023-08-17"
$ws.Range("AB22").Value = "00:00"
$ws.Range("AC22").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD22").Value = $false
$ws.Range("AE22").Value = $false
$ws.Range("AG22").Value = $false
$ws.Range("AW22").Value = "Mimmi Persson"
$ws.Range("AX22").Value = "Mimmi Persson"

# Cells AT22 / AY22 are blank "text" cells in the source diff (no visible
# content). Leave them unset - an empty cell is value-equivalent.

# Restore the default cell style on the cells we temporarily reformatted
# as text, so no stray visual formatting is left behind.
foreach ($addr in @("I22", "Y22", "Z22", "AA22", "AB22")) {
    $ws.Range($addr).Style = "Normal"
}
